$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 262
$ws.Range("F4").Value = 9738
$ws.Range("F5").Value = 663
$ws.Range("F7").Value = 339
$ws.Range("F8").Value = 375
$ws.Range("F11").Value = 206
$ws.Range("F12").Value = 471
$ws.Range("F13").Value = 12358
$ws.Range("F23").Value = 165
$ws.Range("F26").Value = 82
$ws.Range("F28").Value = 60
$ws.Range("F29").Value = 2157
$ws.Range("F30").Value = 1042
$ws.Range("F31").Value = 4219
$ws.Range("F32").Value = 3704
$ws.Range("F33").Value = 666
$ws.Range("F36").Value = 43
$ws.Range("F37").Value = 1337
$ws.Range("F39").Value = 778
$ws.Range("F41").Value = 121
$ws.Range("F42").Value = 449
$ws.Range("F43").Value = 578
$ws.Range("F48").Value = 137
$ws.Range("F49").Value = 152

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 22
$ws.Range("F6").Value = 47
$ws.Range("F14").Value = 40

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 262
$ws.Range("F5").Value = 9738
$ws.Range("F6").Value = 663
$ws.Range("F7").Value = 47
$ws.Range("F9").Value = 339
$ws.Range("F10").Value = 375
$ws.Range("F13").Value = 206
$ws.Range("F14").Value = 471
$ws.Range("F15").Value = 12358
$ws.Range("F23").Value = 165
$ws.Range("F26").Value = 82
$ws.Range("F27").Value = 60
$ws.Range("F28").Value = 2157
$ws.Range("F29").Value = 1042
$ws.Range("F30").Value = 4219
$ws.Range("F31").Value = 3704
$ws.Range("F32").Value = 666
$ws.Range("F35").Value = 43
$ws.Range("F36").Value = 1337
$ws.Range("F38").Value = 778
$ws.Range("F40").Value = 121
$ws.Range("F41").Value = 449
$ws.Range("F43").Value = 578
$ws.Range("F48").Value = 137
$ws.Range("F49").Value = 152
